$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# These cells hold numeric-looking text (e.g. "39.1") stored as shared
# strings, not real numbers. A leading apostrophe forces the new value to
# stay text (matching the original t="s" cell type) instead of being
# auto-coerced to a number; re-applying the "Normal" style afterwards
# strips the quote-prefix formatting flag Excel adds for that apostrophe,
# so the cell's style/number format is left exactly as it was before.

function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.Value = "'" + $val
    $rng.Style = "Normal"
}

# Enterprises density (per 1000 people) - row 12
Set-TextValue "B12" "39.13"
Set-TextValue "C12" "2.64"
Set-TextValue "D12" "41.77"

# Employment (% of total) - row 13
Set-TextValue "D13" "48.83"

# Enterprises (% of total) - row 14
Set-TextValue "B14" "93.27"
Set-TextValue "C14" "6.29"
Set-TextValue "D14" "99.56"
